$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CreateUserSpCharError")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Thu Jul 06 13:04:15 EDT 2023"
